# Auto-update draw results: append the 2025-11-04 Pick 4 draw as a new row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = $ws.UsedRange.Rows.Count + 1

# Assign values with a leading apostrophe so Excel stores them as literal
# text (matching the source data, which keeps date-like / numeric-looking
# strings such as "2025-11-04" and "251104" as text, not as a date serial
# or a number).
$ws.Cells.Item($newRow, 1).Value = "'2025-11-04"
$ws.Cells.Item($newRow, 2).Value = "Pick 4"
$ws.Cells.Item($newRow, 3).Value = "'251104"
$ws.Cells.Item($newRow, 4).Value = "6-9-2-9"
$ws.Cells.Item($newRow, 5).Value = "2025-11-04T21:39:26.734+04:00"

# Reset to the default "Normal" style so the new row does not pick up the
# quote-prefix formatting that typing a leading apostrophe would normally
# apply, keeping it consistent with the rest of the sheet (no explicit
# cell style).
$ws.Range("A" + $newRow + ":E" + $newRow).Style = "Normal"
